$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column for all existing data rows (2-454)
#    from 45190 to 45192.
$ws.Range("C2:C454").Value2 = 45192

# 2. Force row 454 to carry an explicit (custom) row height, matching the
#    target file where row 454 gains ht="15" customHeight="1".
$ws.Rows.Item(454).RowHeight = 15

# 3. Append the new record as row 455.
$ws.Cells.Item(455, 1).Value = "A 44677-2023"
$ws.Cells.Item(455, 2).Value2 = 45189
$ws.Cells.Item(455, 3).Value2 = 45192
$ws.Cells.Item(455, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(455, 5).Value = "ÅSELE"
$ws.Cells.Item(455, 7).Value2 = 7.9
$ws.Cells.Item(455, 8).Value2 = 0
$ws.Cells.Item(455, 9).Value2 = 0
$ws.Cells.Item(455, 10).Value2 = 0
$ws.Cells.Item(455, 11).Value2 = 0
$ws.Cells.Item(455, 12).Value2 = 0
$ws.Cells.Item(455, 13).Value2 = 0
$ws.Cells.Item(455, 14).Value2 = 0
$ws.Cells.Item(455, 15).Value2 = 0
$ws.Cells.Item(455, 16).Value2 = 0
$ws.Cells.Item(455, 17).Value2 = 0

# Apply the date format (style index 1, numFmtId 165 "YYYY-MM-DD") used by
# the other date cells in columns B and C.
$ws.Range("B455:C455").NumberFormat = "YYYY-MM-DD"

# Column R (Artnamn) is left blank but keeps the wrap-text formatting that
# is applied across the whole column.
$ws.Cells.Item(455, 18).WrapText = $true
